# Applies a row-level reshuffle of the weekly price data (rows 2-16).
# Columns D, I, J, K, L, M, N, O, P, Q are permuted across rows according
# to the mapping below (columns A, B, C, E, F, G, H, R are identical on
# every row, so they do not need to move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values currently in the source
# row should end up in the destination row).
$mapping = @{
    2  = 4
    3  = 16
    4  = 15
    5  = 11
    6  = 12
    7  = 5
    8  = 6
    9  = 7
    10 = 2
    11 = 3
    12 = 13
    13 = 9
    14 = 8
    15 = 14
    16 = 10
}

$cols = @("D", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values of every cell that will move,
# so writes to destination rows don't clobber values still needed as a
# source for another destination row.
$snapshot = @{}
foreach ($row in 2..16) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value2 = $snapshot[$srcAddr]
    }
}
